$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.35
$ws.Range("H3").Value = 2.8
$ws.Range("I3").Value = 3.7
$ws.Range("J3").Value = 3.25
$ws.Range("K3").Value = 1.8
$ws.Range("L3").Value = 4.5
$ws.Range("O3").Value = 1.67
$ws.Range("P3").Value = 2.1
$ws.Range("S3").Value = 3.4
$ws.Range("T3").Value = 1.33
$ws.Range("V3").Value = 1.16
$ws.Range("W3").Value = 7
$ws.Range("X3").Value = 1.1
$ws.Range("Y3").Value = 1.73
$ws.Range("Z3").Value = 2
$ws.Range("AA3").Value = 2.5
$ws.Range("AB3").Value = 1.5
$ws.Range("AC3").Value = 5
$ws.Range("AD3").Value = 9
$ws.Range("AI3").Value = 5
$ws.Range("AK3").Value = 23
$ws.Range("AL3").Value = 101
$ws.Range("AP3").Value = 15
$ws.Range("AQ3").Value = 41
$ws.Range("AR3").Value = 41

# Row 4
$ws.Range("G4").Value = 2.45
$ws.Range("I4").Value = 3.5
$ws.Range("J4").Value = 3.4
$ws.Range("Y4").Value = 1.75
$ws.Range("Z4").Value = 2.05
$ws.Range("AA4").Value = 2.5
$ws.Range("AB4").Value = 1.5
$ws.Range("AD4").Value = 9.5
$ws.Range("AF4").Value = 23
$ws.Range("AN4").Value = 7

# Row 15
$ws.Range("H15").Value = 3.9
$ws.Range("I15").Value = 4.5
$ws.Range("O15").Value = 1.14
$ws.Range("P15").Value = 5.5
$ws.Range("S15").Value = 1.53
$ws.Range("U15").Value = 1.79
$ws.Range("V15").Value = 1.94
$ws.Range("Y15").Value = 1.29
$ws.Range("Z15").Value = 3.5
$ws.Range("AA15").Value = 1.53
$ws.Range("AB15").Value = 2.38
$ws.Range("AC15").Value = 10
$ws.Range("AD15").Value = 10
$ws.Range("AI15").Value = 17
$ws.Range("AK15").Value = 13
$ws.Range("AL15").Value = 41
$ws.Range("AM15").Value = 126
$ws.Range("AN15").Value = 17

# Row 16
$ws.Range("G16").Value = 3.6
$ws.Range("L16").Value = 2.63
$ws.Range("N16").Value = 12
$ws.Range("S16").Value = 1.8
$ws.Range("T16").Value = 2
$ws.Range("AH16").Value = 34
$ws.Range("AI16").Value = 12

# Row 18
$ws.Range("G18").Value = 1.3
$ws.Range("H18").Value = 4.75
$ws.Range("I18").Value = 8.75
$ws.Range("J18").Value = 1.78
$ws.Range("K18").Value = 2.37
$ws.Range("L18").Value = 7.5
$ws.Range("O18").Value = 1.22
$ws.Range("P18").Value = 3.45
$ws.Range("S18").Value = 1.65
$ws.Range("T18").Value = 1.98
$ws.Range("W18").Value = 2.55
$ws.Range("X18").Value = 1.39
$ws.Range("AA18").Value = 2.07
$ws.Range("AB18").Value = 1.6
$ws.Range("AC18").Value = 6.6
$ws.Range("AD18").Value = 5.9
$ws.Range("AE18").Value = 8.75
$ws.Range("AF18").Value = 7.7
$ws.Range("AG18").Value = 11.5
$ws.Range("AH18").Value = 32
$ws.Range("AI18").Value = 12
$ws.Range("AJ18").Value = 9.75
$ws.Range("AK18").Value = 24
$ws.Range("AL18").Value = 120
$ws.Range("AN18").Value = 22
$ws.Range("AO18").Value = 60
$ws.Range("AP18").Value = 28
$ws.Range("AQ18").Value = 250
$ws.Range("AS18").Value = 100

Write-Output "Applied all changes"